$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 222, shifting the existing rows 222-281
# (and the sheet dimension) down by one, to make room for a new weekly
# price observation for Espinaca at "Terminal La Palmera de La Serena".
$ws.Rows.Item(222).Insert()

# Fill in the new record (date 2022-06-24, serial 44736).
$ws.Range("A222").Value = 8
$ws.Range("B222").Value = "Terminal La Palmera de La Serena"
$ws.Range("C222").Value = "Coquimbo"
$ws.Range("D222").Value = 44736
$ws.Range("E222").Value = 4
$ws.Range("F222").Value = 100112012
$ws.Range("G222").Value = "Espinaca"
$ws.Range("H222").Value = "Sin especificar"
$ws.Range("I222").Value = "Primera"
$ws.Range("J222").Value = 3000
$ws.Range("K222").Value = 500
$ws.Range("L222").Value = 600
$ws.Range("M222").Value = 550
$ws.Range("N222").Value = "$/atado 300 a 500 gramos"
$ws.Range("O222").Value = "Provincia del Elquí"
$ws.Range("P222").Value = 1100
$ws.Range("Q222").Value = 0.5
$ws.Range("R222").Value = "Hortaliza"
